$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.994.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.649.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.45%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.507'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.30%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("E9").Value = '  +1.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.63%  '
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.887.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.650.53'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.69%  '
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.977.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("E20").Value = '  +1.32%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.42%  '
$ws.Range("E23").Value = '  +4.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.91%  '
$ws.Range("E30").Value = '  +0.59%  '
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.510.58'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.79%  '
$ws.Range("E34").Value = '  +4.52%  '
$ws.Range("E35").Value = '  +8.77%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.573'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.884'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.27%  '
$ws.Range("E39").Value = '  +2.70%  '
$ws.Range("E40").Value = '  +3.53%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("E42").Value = '  +4.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.793.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.775'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.918'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("E51").Value = '  +2.24%  '
